$wb = $excel.ActiveWorkbook

# Rename the "INTER_SWITCH_LINKS" worksheet to "SWITCH_TO_SWITCH"
$switchSheet = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$switchSheet.Name = "SWITCH_TO_SWITCH"

# Make the renamed sheet the active tab (was previously HARDWARE_MANAGEMENT),
# and move its selection to F41
$switchSheet.Activate()
$null = $switchSheet.Range("F41").Select()
